# "Se agrega vista Enviar Acuerdos"
# Update existing applicant row (row 2) with new data, add a second
# applicant row (row 3), re-point the H2 mailto hyperlink to the new
# e-mail address, add a matching hyperlink for the new row, and update
# the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remember the current hyperlink cell style (the "Hipervinculo" style
# already used by H2) so it can be re-applied after adding hyperlinks,
# which otherwise would stamp the cell with a brand new style.
$linkStyle = $ws.Range("H2").Style

# --- Row 2: update existing applicant ---
$ws.Range("C2").Value = 53071568
$ws.Range("D2").Value = "Lili"
$ws.Range("E2").Value = "Llanos"
$ws.Range("G2").Value = 3013109613

$ws.Range("H2").Hyperlinks.Delete()
$ws.Range("H2").Value = "lili@gmail.com"
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:lili@gmail.com") | Out-Null
$ws.Range("H2").Style = $linkStyle

$ws.Range("I2").Value = "Axede"
$ws.Range("J2").Value = 8179562531

# --- Row 3: new applicant ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 9846
$ws.Range("D3").Value = "f"
$ws.Range("E3").Value = "g"
$ws.Range("F3").Value = 2565855
$ws.Range("G3").Value = 301232562

$ws.Range("H3").Value = "h@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:h@hotmail.com") | Out-Null
$ws.Range("H3").Style = $linkStyle

$ws.Range("I3").Value = "s"
$ws.Range("J3").Value = 66616161661

# --- Update the active selection saved with the sheet view ---
$ws.Range("F10").Select() | Out-Null
